$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while forcing text storage
# (preventing Excel's auto date-detection/conversion for strings like "09/04/20")
# and without leaving a residual NumberFormat/style on the cell.
function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2: CODX BUY
Set-TextCell "A2" "CODX"
Set-TextCell "B2" "09/04/20"
Set-TextCell "C2" "BUY"
$ws.Range("D2").Value = 32
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = -384
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 384
$ws.Range("I2").Value = 32
# J2 left blank (no realized profit for a BUY)

# Row 3: CODX SELL
Set-TextCell "A3" "CODX"
Set-TextCell "B3" "09/04/20"
Set-TextCell "C3" "SELL"
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 144
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = -240

# Row 4: Y BUY
Set-TextCell "A4" "Y"
Set-TextCell "B4" "09/04/20"
Set-TextCell "C4" "BUY"
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 123
$ws.Range("F4").Value = -3936
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = 3936
$ws.Range("I4").Value = 32
# J4 left blank (no realized profit for a BUY)

# Row 5: Y SELL
Set-TextCell "A5" "Y"
Set-TextCell "B5" "09/04/20"
Set-TextCell "C5" "SELL"
$ws.Range("D5").Value = 32
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 32
$ws.Range("G5").Value = 122
$ws.Range("H5").Value = 3904
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 0

# Row 6: Y SELL (sell all remaining)
Set-TextCell "A6" "Y"
Set-TextCell "B6" "09/04/20"
Set-TextCell "C6" "SELL"
$ws.Range("D6").Value = 123
$ws.Range("E6").Value = 122
$ws.Range("F6").Value = 15006
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 11102
